# Update the HotStock Top20 list (rows 2-21, columns A-C) to reflect the
# latest rankings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @("长电科技", "长电科技", "金风科技")
    3  = @("特变电工", "特变电工", "海格通信")
    4  = @("三安光电", "通富微电", "岩山科技")
    5  = @("航天发展", "金风科技", "利欧股份")
    6  = @("金风科技", "兆易创新", "特变电工")
    7  = @("岩山科技", "三安光电", "航天发展")
    8  = @("通富微电", "岩山科技", "长电科技")
    9  = @("金太阳",   "金太阳",   "盈新发展")
    10 = @("海格通信", "汉缆股份", "华胜天成")
    11 = @("兆易创新", "蓝色光标", "康强电子")
    12 = @("康强电子", "康强电子", "万向钱潮")
    13 = @("蓝色光标", "森源电气", "蓝色光标")
    14 = @("利欧股份", "东方财富", "三花智控")
    15 = @("华胜天成", "华胜天成", "三安光电")
    16 = @("盈新发展", "利欧股份", "中国西电")
    17 = @("汉缆股份", "海格通信", "通富微电")
    18 = @("三花智控", "贵州茅台", "兆易创新")
    19 = @("信立泰",   "盈新发展", "美年健康")
    20 = @("中国西电", "思源电气", "雷科防务")
    21 = @("天银机电", "航天发展", "平潭发展")
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
}
